# edit.ps1
# Applies the Review_125 -> Review_124 paper-swap edit to the active document.
#
# Strategy: every text chunk in this document is unique, so straightforward
# literal Find/Replace (wdReplaceAll) is used for every 1-for-1 text swap.
# The one structural change -- the old body paragraph shrinks from five
# <w:t> chunks to four -- is handled by locating the chunk that disappears
# and deleting it together with the pair of <w:br/> line breaks that used
# to follow it.

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $ok = $d.Content.Find.Execute(
        $findText,
        $false, $false, $false, $false, $false, $true, 1, $false,
        $replaceText,
        2)
    Write-Host "Replaced? $ok :: $($findText.Substring(0, [Math]::Min(40, $findText.Length)))"
}

# --- 1. Heading paragraph: title + huggingface link -----------------------
Replace-Text `
    "Review 125: Convolutions Die Hard: Open-Vocabulary Segmentation with Single Frozen Convolutional CLIP, 15.08.23" `
    "Review 124: [Short] Self-Alignment with Instruction Backtranslation, 14.08.23"

Replace-Text `
    "https://huggingface.co/papers/2308.02487" `
    "https://huggingface.co/papers/2308.06259"

# --- 2. Bold "Paper:" link --------------------------------------------------
Replace-Text `
    "Paper: https://arxiv.org/abs/2308.02487v2" `
    "Paper: https://arxiv.org/abs/2308.06259v3"

# --- 3. Body paragraph rewrite ---------------------------------------------
# Segment 1
Replace-Text `
    "בד״כ סגמנטציה בתמונות מתבצעת ב-2 שלבים. קודם מחשבים מסכות עבור כל האובייקטים בתמונה ובשלב השני מזהים סוגי האובייקטים. בזמן האחרון יש שימוש רב במודלים מאומנים(כמו CLIP) להפקה של ייצוג התמונה; בשלב 1 מזינים את התמונה למודל המאומן ובשלב 2 מזינים אותה יחד עם המסכות. " `
    "נניח שאתם רוצים לבצע יישור (alignment) עם דאטה מתויג (בסגנון של instruction tuning) של מודל השפה שלכם שאימנתם קודם על דאטה גדול ולא מתויג. נניח שיש בידיכם דאטה מתויג איכותי לא גדול במיוחד ודאטהסט מאוד גדול ומגוון אך לא מתויג. "

# Segment 2
Replace-Text `
    "היום ב #shorthebrewpapereviews סוקרים מאמר המבצע זאת בשלב אחד. למה זה טוב בעצם? כי במקרה הזה צריך להזין את התמונה ל-CLIP רק פעם אחת שזה מקצר משמעותית את זמן ההסקה והאימון כי CLIP זה מודל גדול וכבד. איך הם עשו זאת? " `
    "המאמר שנסקו היום ב-shorthebrewpapereviews מציע שיטה אינטואיטיבית ואלגנטית להפקה של דאטהסט איכותי מתויג בגודל משמעותי מהדאטהסט הלא מתיוג שיש ברשותנו. תהליך האימון מורכב משני שלבים עיקריים: קודם כל מכיילים מודל מאומן ליצור הוראה (instruction) מהתשובה עם הדאטהסט האיכותי המתויג שיש לנו. "

# Segment 3
Replace-Text `
    "קודם כל מעבירים את התמונה דרך CLIP מוקפא ואז מזינים את הייצוג המופק איתו ל-Pixel Decoder יחד עם ״שאילתות האובייקטים״ (סוג של פרומפט לחיפוש האובייקט) לחיזוי המסכות. במהלך האימון מבצעים התאמה בין המסכות ground-truth לבין המסכות שהוצאנו באמצעות אלגוריתם התאמה הונגרי (מזווגים מסכות הדומות ביותר). לאחר שבנינו את המסכות אנו צריכים לזהות את התוכן בתוך המסכות שמצאנו. " `
    "לאחר מכן מזינים למודל את ה״תשובות״ מהדאטהסט הלא מתויג כדי ליצור הוראה לכל לכל תשובה. שלב הזה נקרא self-augmentation. כמובן שלא כל הזוגות שיצרנו הם באיכות גבוהה ואנו מפלטרים אותם בשלב השני הנקרא self-curation. לוקחים מודל שמכויל רק עם הדוגמאות מהדאטהסט המתויג האיכותי (הקטן). מבקשים את המודל (עם פרומפט ספציפי) לדרג מ-1 עד 5 את התאמת התשובה להוראה. "

# Segment 4 (old) is removed entirely, together with the following pair of
# line breaks, since the paragraph shrinks from 5 text chunks to 4.
$rng = $d.Content
$found = $rng.Find.Execute(
    "עבור סגמנטציה עם מילון סגור (הקטגוריות ידועות) מצליבים את הייצוג (שיכון) של הקטגוריה במסכה (המופק באמצעות הפעלת רשת mask pooling) על הפלט של pixel decoder) עם ייצוג הטקטס (המופק עם CLIP) של כל קטגוריה במטרה למצוא קטגוריה הטובה ביותר לכל מסכה. זה נעשה באמצעות חישוב דמיון cosine (עם טמפרטורה נלמדת) בין ייצוגים אלו כאשר קטגוריה עם דמיון מקסימלי עם ייצוג נבחרת כקטגוריה של המסכה. ",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Segment 4 located for deletion? $found"
if ($found) {
    # Extend the matched range by 2 characters so the following
    # <w:br/><w:br/> pair is swallowed along with the text.
    $delRange = $d.Range($rng.Start, $rng.End + 2)
    $delRange.Text = ""
}

# Segment 5 (old) -> becomes the new, final segment of the paragraph
Replace-Text `
    "כדי לאפשר אוגמנטציה עם מילון פתוח (עבור קטגוריות שלא אומנו במהלך האימון) המחברים יוצרים ייצוג המסכה (מוסיפים ״נתיב״ למודל המקביל לנתיב המילון הסגור) רק מהשיכון המופק מהזנת התמונה ל-CLIP (דרך mask pooling). ואז ב-inference משתמשים באותה שיטה שהסברנו עבור המילון הסגור. בשלב משלבים את החיזוי עבור המילון הסגור והפתוח דרך geometric ensemble (סוג של מיצוע)." `
    "לאחר מכן מפלטרים את הזוגות בעלי ציונים הנמוכים. המחברים גם הציעו מה ש נקרא iterative self-curation שבמהלכו לוקחים את הזוגות (הוראה, תשובה) בעלי ציונים גבוהים, מכיילים את המודל עם זה. לאחר מכן ניתן למנף את המודל ל-self-augmentation (השלב הראשון) כדי ליצור דאטהסט מתויג איכותי עוד יותר. ניתן לחזור על התהליך כמה פעמים בתקווה לקבל דאטהסט מתויג גדול ואיכותי. מאמר אלגנטי ונחמד…"

Write-Host "Edit complete."
